# ahmetReflection.docx - "stockIT's Group reflection" cleanup pass.
#
# The only substantive wording fix in this revision lives in the opening
# paragraph: a duplicated "with" and a stray double space are cleaned up.
# The rest of the body text is untouched.

$d = $word.ActiveDocument

# "A team with with initiative" -> "A team with initiative"
$d.Content.Find.Execute("A team with with initiative", $true, $false, $false, $false, $false, $true, 1, $false, "A team with initiative", 2)

# "This  advantage" (double space) -> "This advantage"
$d.Content.Find.Execute("This  advantage", $true, $false, $false, $false, $false, $true, 1, $false, "This advantage", 2)

# Re-touch the section/page setup so Word re-normalises sectPr (this is a
# plain reflow of the section properties - no visual change, matches the
# page geometry the document already had) and records the column spacing
# that Word always emits for a single-column section.
$ps = $d.Sections(1).PageSetup
$ps.SectionStart = 0
$ps.SectionStart = 2
$ps.TextColumns.Spacing = 36
